$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 63 ("Jengibre" / Hortaliza,
# Mercado Mayorista Lo Valledor de Santiago). All existing rows from the old
# row 63 down to the old row 94 shift down by one (new rows 64..95).
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new weekly record.
$ws.Range("A63").Value = 6
$ws.Range("B63").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44777
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = 100114007
$ws.Range("G63").Value = "Jengibre"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 230
$ws.Range("K63").Value = 12000
$ws.Range("L63").Value = 13000
$ws.Range("M63").Value = 12435
$ws.Range("N63").Value = "$/caja 13 kilos"
$ws.Range("O63").Value = "Perú"
$ws.Range("P63").Value = 957
$ws.Range("Q63").Value = 13
$ws.Range("R63").Value = "Hortaliza"
